# Update "想去人数" (interest count) figures in column F across all four
# sheets of the 上海-漫展信息 workbook, matching the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2106
$ws1.Range("F5").Value = 805
$ws1.Range("F6").Value = 40573
$ws1.Range("F7").Value = 1439
$ws1.Range("F10").Value = 874
$ws1.Range("F11").Value = 5574
$ws1.Range("F12").Value = 386
$ws1.Range("F14").Value = 2679
$ws1.Range("F15").Value = 6282
$ws1.Range("F17").Value = 1172
$ws1.Range("F18").Value = 677
$ws1.Range("F21").Value = 1086
$ws1.Range("F27").Value = 847
$ws1.Range("F33").Value = 13
$ws1.Range("F34").Value = 178
$ws1.Range("F36").Value = 188
$ws1.Range("F37").Value = 1138
$ws1.Range("F38").Value = 36
$ws1.Range("F39").Value = 64

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 501
$ws2.Range("F21").Value = 124
$ws2.Range("F28").Value = 434
$ws2.Range("F29").Value = 928
$ws2.Range("F30").Value = 539
$ws2.Range("F32").Value = 69
$ws2.Range("F36").Value = 114

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 811
$ws3.Range("F6").Value = 506

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 811
$ws4.Range("F8").Value = 506
$ws4.Range("F11").Value = 501
$ws4.Range("F12").Value = 805
$ws4.Range("F13").Value = 1439
$ws4.Range("F16").Value = 5574
$ws4.Range("F17").Value = 386
$ws4.Range("F19").Value = 2679
$ws4.Range("F21").Value = 6282
$ws4.Range("F24").Value = 1172
$ws4.Range("F27").Value = 677
$ws4.Range("F29").Value = 1086
$ws4.Range("F30").Value = 124
$ws4.Range("F33").Value = 847
$ws4.Range("F36").Value = 1113
$ws4.Range("F39").Value = 928
$ws4.Range("F40").Value = 539
$ws4.Range("F41").Value = 178
$ws4.Range("F43").Value = 69
$ws4.Range("F44").Value = 188
$ws4.Range("F47").Value = 0
$ws4.Range("F49").Value = 64
